# Auto-update predictions and index for 2025-10-22
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# NOTE: the "Win %" column stores its values as literal text (e.g. "74%"),
# not real percentage numbers. Assigning a bare "NN%" string lets Excel's
# type inference convert it to a numeric percentage, so each of those
# cells is temporarily switched to a text number format before the
# assignment (forcing a text literal, matching the source data) and then
# restored to the Normal style so no stray formatting is left behind.
function Set-TextPercent($cell, $text) {
    $cell.NumberFormat = "@"
    $cell.Value = $text
    $cell.Style = "Normal"
}

# Row 2: Bayern Munich - Club Brugge KV
$ws.Range("B2").Value = "Bayern Munich  - Club Brugge KV: 20:00"
$ws.Range("C2").Value = 2.6
$ws.Range("D2").Value = "Bayern Munich"
Set-TextPercent $ws.Range("F2") "74%"
$ws.Range("H2").Value = 20

# Row 3: CF América - Puebla FC (final score, checkmark result)
$ws.Range("B3").Value = "CF América ✓ - Puebla FC: 2:1"
$ws.Range("C3").Value = 2.22
$ws.Range("D3").Value = "CF América"
Set-TextPercent $ws.Range("F3") "73%"
$ws.Range("G3").Value = "✓"
$ws.Range("H3").Value = 3
$ws.Range("I3").Value = $true

# Row 4: Real Madrid - Juventus FC
$ws.Range("B4").Value = "Real Madrid  - Juventus FC: 20:00"
$ws.Range("C4").Value = 3.04
$ws.Range("D4").Value = "Real Madrid"
$ws.Range("E4").Value = 4.5
Set-TextPercent $ws.Range("F4") "72%"
$ws.Range("H4").Value = 20
$ws.Range("I4").Value = $false

# Row 5: Chelsea FC - Ajax Amsterdam
$ws.Range("B5").Value = "Chelsea FC  - Ajax Amsterdam: 20:00"
$ws.Range("C5").Value = 2.03
$ws.Range("D5").Value = "Chelsea FC"
Set-TextPercent $ws.Range("F5") "70%"
$ws.Range("H5").Value = 20

# Row 6: FC Santa Coloma - Penya Encarnada d'Andorra (final score, checkmark result)
$ws.Range("B6").Value = "FC Santa Coloma ✓ - Penya Encarnada d'Andorra: 1:0"
$ws.Range("C6").Value = 2.4
$ws.Range("D6").Value = "FC Santa Coloma"
$ws.Range("E6").Value = 3.5
Set-TextPercent $ws.Range("F6") "60%"
$ws.Range("G6").Value = "✓"
$ws.Range("H6").Value = 1
$ws.Range("I6").Value = $true

# Row 7: Nagaworld FC - Phnom Penh Crown
$ws.Range("B7").Value = "Nagaworld FC - Phnom Penh Crown : 12:00"
$ws.Range("C7").Value = 2.21
$ws.Range("D7").Value = "Phnom Penh Crown"
Set-TextPercent $ws.Range("F7") "60%"
$ws.Range("H7").Value = 12

# Row 8: JS Kabylie - MC El Bayadh
$ws.Range("B8").Value = "JS Kabylie  - MC El Bayadh: 23:00"
$ws.Range("C8").Value = 0.88
$ws.Range("D8").Value = "JS Kabylie"
$ws.Range("E8").Value = 1.5
Set-TextPercent $ws.Range("F8") "58%"
$ws.Range("H8").Value = 23

# Row 9: MC Algiers - JS Saoura
$ws.Range("B9").Value = "MC Algiers  - JS Saoura: 23:00"
$ws.Range("C9").Value = 1.07
$ws.Range("D9").Value = "MC Algiers"
$ws.Range("E9").Value = 2.5
Set-TextPercent $ws.Range("F9") "55%"
$ws.Range("H9").Value = 23

# Row 10: Deportivo La Guaira - Zamora FC (final score)
$ws.Range("B10").Value = "Deportivo La Guaira  - Zamora FC: 1:1"
$ws.Range("C10").Value = 1.4
$ws.Range("D10").Value = "Deportivo La Guaira"
$ws.Range("E10").Value = 2.5
Set-TextPercent $ws.Range("F10") "53%"
$ws.Range("H10").Value = 2
$ws.Range("I10").Value = $true

# Row 11: MC Algiers - Paradou AC
$ws.Range("B11").Value = "MC Algiers  - Paradou AC: 20:00"
$ws.Range("C11").Value = 1.53
$ws.Range("D11").Value = "MC Algiers"
$ws.Range("E11").Value = 2.5
Set-TextPercent $ws.Range("F11") "52%"
$ws.Range("H11").Value = 20

# Drop the now-stale fixture rows (old rows 12-17); this shifts the
# summary formulas in rows 18-20 up to rows 12-14 and Excel auto-adjusts
# the K18/K20 style references inside them.
$ws.Rows("12:17").Delete()
